$d = $word.ActiveDocument

function Find-ParagraphIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match $pattern) {
            return $i
        }
    }
    return -1
}

# 1) After the "14.Social Link for Contact" paragraph, add a new paragraph
#    containing "15.Quick Link" (inheriting the Cambria/Bold formatting),
#    followed by a new blank paragraph with the same formatting.
$idx14 = Find-ParagraphIndex("14\.Social Link for Contact")
if ($idx14 -gt 0) {
    $p14 = $d.Paragraphs.Item($idx14)
    $p14.Range.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($idx14 + 1)
    $newP.Range.InsertBefore("15.Quick Link")
}

# 2) Paragraph "2. Register " -> "2. Logout "
$idxReg = Find-ParagraphIndex("^2\. Register ")
if ($idxReg -gt 0) {
    $p = $d.Paragraphs.Item($idxReg)
    $full = $p.Range.Text
    $offset = $full.IndexOf(" Register ")
    $start = $p.Range.Start + $offset
    $end = $start + " Register ".Length
    $target = $d.Range($start, $end)
    $target.Text = " Logout "
}

# 3) Paragraph "9. Zip " -> "9. Thank you Message "
$idxZip = Find-ParagraphIndex("^9\. Zip ")
if ($idxZip -gt 0) {
    $p = $d.Paragraphs.Item($idxZip)
    $full = $p.Range.Text
    $offset = $full.IndexOf("Zip ")
    $start = $p.Range.Start + $offset
    $end = $start + "Zip ".Length
    $target = $d.Range($start, $end)
    $target.Text = "Thank you Message "
}

# 4) Paragraph "11. Account number" -> "11. Product Image for checkout. "
$idxAcct = Find-ParagraphIndex("^11\. Account number")
if ($idxAcct -gt 0) {
    $p = $d.Paragraphs.Item($idxAcct)
    $full = $p.Range.Text
    $offset = $full.IndexOf("Account number")
    $start = $p.Range.Start + $offset
    $end = $start + "Account number".Length
    $target = $d.Range($start, $end)
    $target.Text = "Product Image for checkout. "
}
